# Updates the cryptocurrency price ("Price", column D) and 1h volume change
# ("Volume(1h)", column E) figures on the active worksheet to the refreshed
# values from the latest data pull.
#
# Column D cells that are not purely numeric-looking text (e.g. thousand-dot
# separated prices like "43.617.17") are written as-is; cells whose new
# value would otherwise be auto-recognized as a number by Excel (e.g.
# "309.68") are entered with a leading apostrophe so they stay plain text,
# just like typing them in by hand.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "43.617.17"
$ws.Cells.Item(2, 5).Value = "  +1.02%  "
$ws.Cells.Item(3, 4).Value = "2.376.43"
$ws.Cells.Item(3, 5).Value = "  +2.95%  "
$ws.Cells.Item(4, 5).Value = "  -0.01%  "
$ws.Cells.Item(5, 4).Value = "'309.68"
$ws.Cells.Item(5, 5).Value = "  -0.30%  "
$ws.Cells.Item(6, 4).Value = "'104.69"
$ws.Cells.Item(6, 5).Value = "  +3.23%  "
$ws.Cells.Item(7, 4).Value = "'0.508"
$ws.Cells.Item(7, 5).Value = "  -5.53%  "
$ws.Cells.Item(8, 5).Value = "  -0.02%  "
$ws.Cells.Item(9, 5).Value = "  -1.02%  "
$ws.Cells.Item(10, 4).Value = "'36.03"
$ws.Cells.Item(10, 5).Value = "  +0.12%  "
$ws.Cells.Item(11, 4).Value = "'53.42"
$ws.Cells.Item(11, 5).Value = "  +2.18%  "
$ws.Cells.Item(12, 4).Value = "'0.0813"
$ws.Cells.Item(12, 5).Value = "  -0.71%  "
$ws.Cells.Item(13, 5).Value = "  -0.79%  "
$ws.Cells.Item(14, 5).Value = "  -2.36%  "
$ws.Cells.Item(15, 4).Value = "2.742.54"
$ws.Cells.Item(15, 5).Value = "  +2.96%  "
$ws.Cells.Item(16, 4).Value = "'15.64"
$ws.Cells.Item(16, 5).Value = "  +4.21%  "
$ws.Cells.Item(17, 4).Value = "2.374.36"
$ws.Cells.Item(17, 5).Value = "  +2.89%  "
$ws.Cells.Item(18, 4).Value = "'0.811"
$ws.Cells.Item(18, 5).Value = "  -0.25%  "
$ws.Cells.Item(19, 4).Value = "43.564.35"
$ws.Cells.Item(19, 5).Value = "  +1.07%  "
$ws.Cells.Item(20, 4).Value = "'6.31"
$ws.Cells.Item(20, 5).Value = "  +3.55%  "
$ws.Cells.Item(21, 4).Value = "'11.90"
$ws.Cells.Item(21, 5).Value = "  -5.25%  "
$ws.Cells.Item(22, 4).Value = "0.0₃0917"
$ws.Cells.Item(22, 5).Value = "  -0.48%  "
$ws.Cells.Item(23, 4).Value = "'68.43"
$ws.Cells.Item(23, 5).Value = "  -0.24%  "
$ws.Cells.Item(24, 4).Value = "'240.96"
$ws.Cells.Item(24, 5).Value = "  -0.12%  "
$ws.Cells.Item(25, 5).Value = "  +1.89%  "
$ws.Cells.Item(26, 4).Value = "'2.62"
$ws.Cells.Item(26, 5).Value = "  -0.40%  "
$ws.Cells.Item(27, 5).Value = "  +0.09%  "
$ws.Cells.Item(28, 4).Value = "'25.80"
$ws.Cells.Item(28, 5).Value = "  +3.91%  "
$ws.Cells.Item(29, 5).Value = "  -2.91%  "
$ws.Cells.Item(30, 4).Value = "'36.61"
$ws.Cells.Item(30, 5).Value = "  -2.61%  "
$ws.Cells.Item(31, 5).Value = "  -1.27%  "
$ws.Cells.Item(32, 5).Value = "  -0.28%  "
$ws.Cells.Item(33, 4).Value = "'160.87"
$ws.Cells.Item(33, 5).Value = "  -4.01%  "
$ws.Cells.Item(34, 4).Value = "'5.26"
$ws.Cells.Item(34, 5).Value = "  -1.28%  "
$ws.Cells.Item(35, 4).Value = "'0.999"
$ws.Cells.Item(35, 5).Value = "  -0.03%  "
$ws.Cells.Item(36, 4).Value = "'18.29"
$ws.Cells.Item(36, 5).Value = "  +3.12%  "
$ws.Cells.Item(37, 5).Value = "  +5.73%  "
$ws.Cells.Item(38, 4).Value = "'3.12"
$ws.Cells.Item(38, 5).Value = "  -0.36%  "
$ws.Cells.Item(39, 4).Value = "'0.0740"
$ws.Cells.Item(39, 5).Value = "  -0.25%  "
$ws.Cells.Item(40, 4).Value = "'4.65"
$ws.Cells.Item(40, 5).Value = "  +7.74%  "
$ws.Cells.Item(41, 4).Value = "'1.94"
$ws.Cells.Item(41, 5).Value = "  +5.48%  "
$ws.Cells.Item(42, 5).Value = "  -2.04%  "
$ws.Cells.Item(43, 5).Value = "  -1.99%  "
$ws.Cells.Item(44, 4).Value = "'2.60"
$ws.Cells.Item(44, 5).Value = "  +13.43%  "
$ws.Cells.Item(45, 4).Value = "2.032.54"
$ws.Cells.Item(45, 5).Value = "  +2.33%  "
$ws.Cells.Item(46, 4).Value = "'19.75"
$ws.Cells.Item(46, 5).Value = "  +3.45%  "
$ws.Cells.Item(47, 5).Value = "  +0.16%  "
$ws.Cells.Item(48, 5).Value = "  +3.43%  "
$ws.Cells.Item(49, 4).Value = "'10.57"
$ws.Cells.Item(49, 5).Value = "  +7.26%  "
$ws.Cells.Item(50, 4).Value = "'58.12"
$ws.Cells.Item(50, 5).Value = "  +4.39%  "
$ws.Cells.Item(51, 5).Value = "  +0.32%  "
